$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- B column: best_params ----
$ws.Range("B2").Value = "{'max_depth': 5, 'min_samples_leaf': 2, 'min_samples_split': 10}"
$ws.Range("B3").Value = "{'max_depth': 25, 'min_samples_split': 2, 'n_estimators': 250}"
$ws.Range("B4").Value = "{'learning_rate': 0.05, 'n_estimators': 100, 'num_leaves': 50}"
$ws.Range("B5").Value = "{'learning_rate': 0.1, 'max_depth': 3, 'n_estimators': 150}"

# ---- D column: best_model ----
$ws.Range("D2").Value = "DecisionTreeClassifier(max_depth=5, min_samples_leaf=2, min_samples_split=10)"
$ws.Range("D3").Value = "RandomForestClassifier(max_depth=25, n_estimators=250)"
$ws.Range("D5").Value = "XGBClassifier(base_score=None, booster=None, callbacks=None,`n              colsample_bylevel=None, colsample_bynode=None,`n              colsample_bytree=None, device=None, early_stopping_rounds=None,`n              enable_categorical=True, eval_metric=None, feature_types=None,`n              gamma=None, grow_policy=None, importance_type=None,`n              interaction_constraints=None, learning_rate=0.1, max_bin=None,`n              max_cat_threshold=None, max_cat_to_onehot=None,`n              max_delta_step=None, max_depth=3, max_leaves=None,`n              min_child_weight=None, missing=nan, monotone_constraints=None,`n              multi_strategy=None, n_estimators=150, n_jobs=None,`n              num_parallel_tree=None, random_state=None, ...)"

# ---- E column: confusion_matrix ----
$ws.Range("E3").Value = "[[215  62]`n [ 24 425]]"
$ws.Range("E5").Value = "[[221  56]`n [ 28 421]]"

# ---- Row 2 (CART) ----
$ws.Range("M2").Value = 17.63518118858337

# ---- Row 3 (Random Forest) ----
$ws.Range("C3").Value = 0.8715535015997155
$ws.Range("G3").Value = 62
$ws.Range("I3").Value = 215
$ws.Range("J3").Value = 0.882950251734535
$ws.Range("K3").Value = 0.8815426997245179
$ws.Range("L3").Value = 0.8795854818582091
$ws.Range("M3").Value = 626.4200420379639

# ---- Row 4 (LightGBM) ----
$ws.Range("M4").Value = 217.9192698001862

# ---- Row 5 (XGBoost) ----
$ws.Range("C5").Value = 0.8749946676146463
$ws.Range("F5").Value = 421
$ws.Range("G5").Value = 56
$ws.Range("H5").Value = 28
$ws.Range("I5").Value = 221
$ws.Range("J5").Value = 0.8844884536673994
$ws.Range("K5").Value = 0.8842975206611571
$ws.Range("L5").Value = 0.8829672685356211
$ws.Range("M5").Value = 265.2381844520569

# Entering multi-line text (embedded newlines) causes Excel to auto-expand
# the row height; the source workbook keeps default row heights, so
# re-autofit the affected rows back down.
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(5).AutoFit()
